$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "E7"  = 16.0416
    "B9"  = 5.801999999999995
    "E12" = 18.24390000000003
    "E14" = 16.6393
    "B18" = 6.402399999999997
    "B20" = 9.072400000000002
    "E26" = 15.9598
    "B27" = 6.258200000000004
    "E27" = 16.66929999999999
    "E29" = 17.04130000000002
    "B35" = 8.561400000000004
    "E37" = 16.7676
    "E38" = 16.5028
    "E51" = 17.2009
    "E52" = 16.98220000000001
    "E55" = 16.54720000000001
    "B69" = 5.296599999999995
    "E69" = 17.26050000000002
    "E70" = 18.19980000000002
    "B76" = 5.716900000000002
    "B78" = 9.905499999999998
    "E81" = 16.60859999999999
    "B82" = 5.404100000000001
    "B83" = 5.674700000000001
    "E83" = 16.5792
    "B93" = 6.006799999999997
    "E102" = 16.9058
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
